$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price/volume snapshot (GitHub Actions data refresh).
# Price cells that look like plain numbers ("1.002", "0.9978", ...) are
# forced back to literal text (leading apostrophe) so Excel doesn't
# reinterpret them as numeric values, then the cell style is reset to
# "Normal" so no extra number-format/quote-prefix styling is left behind.
# Two rows (40/41) also swap their Coin/Link content (FraxShare <-> TrustWalletToken).

$ws.Range("D2").Value = "24.397.41"
$ws.Range("E2").Value = "  -4.41%  "
$ws.Range("D3").Value = "1.641.20"
$ws.Range("E3").Value = "  -6.55%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'0.9978"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "'305.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("D7").Value = "'0.3616"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.58%  "
$ws.Range("D8").Value = "'47.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.84%  "
$ws.Range("D9").Value = "'0.3248"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.67%  "
$ws.Range("D10").Value = "'1.114"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.98%  "
$ws.Range("D11").Value = "'0.06897"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -10.12%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "'5.914"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.97%  "
$ws.Range("D14").Value = "'19.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -11.47%  "
$ws.Range("D15").Value = "1.644.02"
$ws.Range("E15").Value = "  -6.28%  "
$ws.Range("D16").Value = "'6.516"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.56%  "
$ws.Range("D17").Value = "'0.00001042"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -9.80%  "
$ws.Range("D18").Value = "'0.06479"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.58%  "
$ws.Range("D19").Value = "'0.9987"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "'76.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -11.18%  "
$ws.Range("D21").Value = "'5.876"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.06%  "
$ws.Range("D22").Value = "'15.66"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -11.30%  "
$ws.Range("D23").Value = "'12.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.21%  "
$ws.Range("D24").Value = "24.390.89"
$ws.Range("E24").Value = "  -4.45%  "
$ws.Range("D25").Value = "'2.400"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("D26").Value = "'2.325"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -19.97%  "
$ws.Range("D27").Value = "'145.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.30%  "
$ws.Range("D28").Value = "'18.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.48%  "
$ws.Range("D29").Value = "1.826.81"
$ws.Range("E29").Value = "  -6.25%  "
$ws.Range("D30").Value = "'123.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.66%  "
$ws.Range("D31").Value = "'1.143"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.32%  "
$ws.Range("D32").Value = "'4.060"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("D33").Value = "'5.584"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -22.38%  "
$ws.Range("D34").Value = "'0.08312"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.13%  "
$ws.Range("D35").Value = "'1.673"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.15%  "
$ws.Range("D36").Value = "'12.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -14.12%  "
$ws.Range("D37").Value = "'5.125"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -11.40%  "
$ws.Range("D38").Value = "'0.06016"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.50%  "
$ws.Range("D39").Value = "'0.02207"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.63%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.193"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -12.66%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.197"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.39%  "
$ws.Range("D42").Value = "'0.2030"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.41%  "
$ws.Range("D43").Value = "'0.9978"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.25%  "
$ws.Range("D44").Value = "'0.5822"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -11.73%  "
$ws.Range("D45").Value = "'3.715"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.86%  "
$ws.Range("D46").Value = "'12.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -12.72%  "
$ws.Range("D47").Value = "'0.5567"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.88%  "
$ws.Range("D48").Value = "'121.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.64%  "
$ws.Range("D49").Value = "'1.925"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.98%  "
$ws.Range("E50").Value = "  -7.72%  "
$ws.Range("D51").Value = "'73.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.50%  "
